# Applies the "favorability_check() only needs 'df' and 'tb' as input" edit:
#  - drops the now-unused nettax_abg_nokfb_tu / nettax_abg_kfb_tu columns (old V/W),
#    collapsing U into a simple "=H<row>" (nettax_kfb_tu) column
#  - simplifies/rewrites a couple of the T/U/O formulas
#  - appends three new sample rows (10-12)
#  - updates the view selection and the stray absPath metadata

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop columns V and W (nettax_kfb_tu / nettax_abg_kfb_tu) entirely - their
#    header strings disappear from the shared-string table once unreferenced.
# ---------------------------------------------------------------------------
$ws.Range("V1:W9").Clear()

# ---------------------------------------------------------------------------
# 2. Column U becomes the (former W-derived) "nettax_kfb_tu" column.
#    Row 2 keeps a plain (non-shared) formula; rows 3-9 share one formula.
# ---------------------------------------------------------------------------
$ws.Range("U1").Value = "nettax_kfb_tu"
$ws.Range("U2").Formula = "=H2"
$ws.Range("U3:U9").Formula = "=H3"

# ---------------------------------------------------------------------------
# 3. O2:O4 now MIN() only over T:U (since V:W are gone).
# ---------------------------------------------------------------------------
$ws.Range("O2").Formula = "=MIN(T2:U2)/12"
$ws.Range("O3").Formula = "=MIN(T3:U3)/12"
$ws.Range("O4").Formula = "=MIN(T4:U4)/12"

# ---------------------------------------------------------------------------
# 4. Row 5 gets bespoke formulas (no longer a MIN of T:W).
# ---------------------------------------------------------------------------
$ws.Range("O5").Formula = "=F5/12"
$ws.Range("T5").Formula = "=F5-12*L5+8000"
$ws.Range("U5").Formula = "=H5+8000"

# ---------------------------------------------------------------------------
# 5. C3:C8 shared "increment" formula (si renumbers naturally once V/W are
#    gone - content/values are what matter).
# ---------------------------------------------------------------------------
$ws.Range("C3:C8").Formula = "=C2+1"

# ---------------------------------------------------------------------------
# 6. Three brand-new sample rows (10-12), cloning formatting from row 9
#    (column-block at a time, so no stray blank N/S cells are materialised)
#    then overwriting values/formulas. Column E loses row 9's TRUE/FALSE
#    number format since the new rows use the plain default style.
# ---------------------------------------------------------------------------
foreach ($r in 10, 11, 12) {
    $ws.Range("A9:M9").Copy()
    $ws.Range("A$r`:M$r").PasteSpecial(-4122) # xlPasteFormats
    $ws.Range("O9:R9").Copy()
    $ws.Range("O$r`:R$r").PasteSpecial(-4122)
    $ws.Range("T9:U9").Copy()
    $ws.Range("T$r`:U$r").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false
    $ws.Range("E$r").ClearFormats()
}

# Row 10
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = $true
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = 12000
$ws.Range("G10").Value = 12000
$ws.Range("H10").Value = 10000
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 194
$ws.Range("M10").Value = 2019
$ws.Range("O10").Formula = "=2*T10/12"
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("T10").Formula = "=F10-(12*L10)"
$ws.Range("U10").Formula = "=H10"

# Row 11
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = 12000
$ws.Range("G11").Value = 12000
$ws.Range("H11").Value = 10000
$ws.Range("I11").Value = 10000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 194
$ws.Range("M11").Value = 2019
$ws.Range("O11").Formula = "=2*T11/12"
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("T11:T12").Formula = "=F11-(12*L11)"
$ws.Range("U11:U12").Formula = "=H11"

# Row 12
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 11
$ws.Range("D12").Value = $false
$ws.Range("E12").Value = $true
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 194
$ws.Range("L12").Value = 194
$ws.Range("M12").Value = 2019
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 0
$ws.Range("T12").Value = 0
# U12 formula already applied above via the U11:U12 shared fill.

# ---------------------------------------------------------------------------
# 7. View: scroll so column C is leftmost and select O5.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("O5").Select()
